$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "E2" = 58
    "E3" = 24
    "E4" = 25
    "E5" = 63
    "E7" = 14
    "E9" = 4
    "E10" = 192
    "F10" = 78
    "H10" = 78
    "E11" = 139
    "F11" = 67
    "H11" = 67
    "E12" = 203
    "F12" = 102
    "H12" = 102
    "E13" = 68
    "F13" = 28
    "H13" = 28
    "E14" = 60
    "F14" = 24
    "H14" = 24
    "E15" = 87
    "E16" = 81
    "F16" = 35
    "H16" = 35
    "E17" = 38
    "E20" = 49
    "E21" = 65
    "E22" = 82
    "F22" = 34
    "H22" = 34
    "E23" = 93
    "F23" = 35
    "H23" = 35
    "E24" = 95
    "F24" = 43
    "H24" = 43
    "E25" = 86
    "F25" = 32
    "H25" = 32
    "E26" = 49
    "F26" = 21
    "H26" = 21
    "E27" = 131
    "F27" = 61
    "H27" = 61
    "E28" = 85
    "F28" = 21
    "H28" = 21
    "E29" = 80
    "E30" = 92
    "F30" = 42
    "H30" = 42
    "E31" = 39
    "F31" = 17
    "H31" = 17
    "E32" = 89
    "F32" = 44
    "H32" = 44
    "E33" = 117
    "E34" = 94
    "E35" = 60
    "E36" = 31
    "E37" = 63
    "F37" = 26
    "H37" = 26
    "F38" = 24
    "H38" = 24
    "E39" = 99
    "E40" = 128
    "F40" = 48
    "H40" = 48
    "E41" = 164
    "F41" = 56
    "H41" = 56
    "E42" = 144
    "F42" = 65
    "H42" = 65
    "E43" = 42
    "E44" = 129
    "F44" = 56
    "H44" = 56
    "E46" = 111
    "F46" = 46
    "H46" = 46
    "E47" = 189
    "F47" = 72
    "H47" = 72
    "E48" = 94
    "E49" = 112
    "F49" = 45
    "H49" = 45
    "E50" = 87
    "F50" = 33
    "H50" = 33
    "E51" = 91
    "F51" = 34
    "H51" = 34
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

$wb.Save()
